# "Changed valve state reading functionality"
#
# - Row 13: the End_Time reading was re-taken (valve was read later) and the
#   End_Volume (mL) syringe reading changed, which ripples through the
#   computed Duration / Volume-by-Syringe / Accuracy columns. The Accuracy
#   formula's constant (tube volume offset) also changed from 141 to 209.
# - Row 14: a brand new valve reading was recorded (previously an empty
#   template row) together with its computed columns, extending the table
#   and the V-column color-scale conditional formatting down to row 14.
# - The active selection/view is reset to the newly-entered cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13: update the re-read End_Time / End_Volume values -------------
$ws.Range("H13").Value = 45870.604166666664
$ws.Range("J13").Value = 6
$ws.Range("U13").Formula = "=Q13*1440/209"

# --- Row 14: new valve reading --------------------------------------------
$ws.Range("B14").Value = 11
$ws.Range("C14").Value = 290
$ws.Range("D14").Value = 163
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 45870.60833333333
$ws.Range("H14").Value = 45870.703472222223
$ws.Range("I14").Value = 6
$ws.Range("J14").Value = 5.48
$ws.Range("K14").Value = 5

$ws.Range("O14").Formula = "=H14-G14"
$ws.Range("P14").Formula = "=O14"
$ws.Range("Q14").Formula = "=I14-J14"
$ws.Range("R14").Formula = "=(F14-E14)/0.9982"
$ws.Range("U14").Formula = "=Q14*1440/137"
$ws.Range("V14").Formula = "=(1-ABS(U14-K14)/K14)*100%"

# --- Extend the Accuracy color-scale conditional formatting to row 14 ----
$ws.Range("V4:V13").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("V4:V14"))

# --- Reset the view / selection to the newly added row -------------------
$ws.Range("R15").Select()
